# Updated cryptos list (price + 1h volume change columns) per the
# Thu Apr 20 21:35:21 UTC 2023 GitHub Actions refresh commit.
#
# Each data row lives in the "Coin" table on the active sheet: column D
# holds the latest price (text, not a number, because coinranking
# formats some prices like "28.497.56" with thousands separators that
# are not valid numeric literals) and column E holds the 1h volume
# change percentage (also stored as text, with padding spaces).
#
# Numeric-looking price strings (e.g. "1.010", "0.4764") would silently
# be re-interpreted by Excel as numbers (dropping the trailing zero /
# the fixed decimal formatting) if we just assign .Value, so for those
# cells we first force the cell to Text format ("@") and only then
# write the string — this keeps the value stored exactly as text,
# matching the original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.457.46"
$ws.Range("E2").Value = "  -3.31%  "
$ws.Range("D3").Value = "1.957.84"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -0.73%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.70"
$ws.Range("E5").Value = "  -1.94%  "
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4764"
$ws.Range("E7").Value = "  -4.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4062"
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.47"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08456"
$ws.Range("E10").Value = "  -4.86%  "
$ws.Range("E11").Value = "  -4.21%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.14"
$ws.Range("E12").Value = "  -4.14%  "
$ws.Range("D13").Value = "1.956.63"
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.605"
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.157"
$ws.Range("E15").Value = "  -3.87%  "
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.45"
$ws.Range("E17").Value = "  -4.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001071"
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06610"
$ws.Range("E19").Value = "  -0.68%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.63"
$ws.Range("E20").Value = "  -3.64%  "
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.825"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").Value = "28.488.55"
$ws.Range("E23").Value = "  -3.49%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.57"
$ws.Range("E24").Value = "  -2.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.288"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").Value = "2.209.90"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.96"
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.21"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.919"
$ws.Range("E29").Value = "  -4.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.157"
$ws.Range("E30").Value = "  -5.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.45"
$ws.Range("E31").Value = "  -2.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9796"
$ws.Range("E32").Value = "  -6.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09603"
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.443"
$ws.Range("E34").Value = "  -6.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.589"
$ws.Range("E35").Value = "  -3.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.659"
$ws.Range("E36").Value = "  -3.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02333"
$ws.Range("E37").Value = "  -4.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.835"
$ws.Range("E38").Value = "  -3.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06214"
$ws.Range("E39").Value = "  -1.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.254"
$ws.Range("E40").Value = "  -3.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6221"
$ws.Range("E41").Value = "  -4.08%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.13"
$ws.Range("E42").Value = "  -3.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.010"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1918"
$ws.Range("E44").Value = "  -5.33%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.339"
$ws.Range("E45").Value = "  +2.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5961"
$ws.Range("E46").Value = "  -4.77%  "
$ws.Range("E47").Value = "  -2.29%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.057"
$ws.Range("E48").Value = "  -5.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.402"
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000331"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06826"
$ws.Range("E51").Value = "  -1.87%  "
